$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# Row 7 (ALC), Leve Item ID 1960
$ws_ALC.Range("H7").Value = 0
$ws_ALC.Range("I7").Value = 0
$ws_ALC.Range("J7").Value = 0
$ws_ALC.Range("K7").Value = 0
$ws_ALC.Range("L7").Value = 0
$ws_ALC.Range("N7").ClearContents()

# Row 14 (ALC), Leve Item ID 1960
$ws_ALC.Range("H14").Value = 0
$ws_ALC.Range("I14").Value = 0
$ws_ALC.Range("J14").Value = 0
$ws_ALC.Range("K14").Value = 0
$ws_ALC.Range("L14").Value = 0
$ws_ALC.Range("N14").ClearContents()

# Row 15 (ALC), Leve Item ID 44146
$ws_ALC.Range("H15").Value = 517.8889
$ws_ALC.Range("I15").Value = 517.8889
$ws_ALC.Range("J15").Value = 0
$ws_ALC.Range("K15").Value = 1553.6667
$ws_ALC.Range("L15").Value = 0
$ws_ALC.Range("M15").Value = -1384.6667

# Row 18 (ALC), Leve Item ID 5471
$ws_ALC.Range("H18").Value = 997
$ws_ALC.Range("I18").Value = 997
$ws_ALC.Range("J18").Value = 0
$ws_ALC.Range("K18").Value = 997
$ws_ALC.Range("L18").Value = 0
$ws_ALC.Range("M18").Value = -713
$ws_ALC.Range("N18").ClearContents()

# Row 21 (ALC), Leve Item ID 2149
$ws_ALC.Range("H21").Value = 3017
$ws_ALC.Range("I21").Value = 3017
$ws_ALC.Range("J21").Value = 0
$ws_ALC.Range("K21").Value = 3017
$ws_ALC.Range("L21").Value = 0
$ws_ALC.Range("M21").Value = -2549

# Row 23 (ALC), Leve Item ID 2149
$ws_ALC.Range("H23").Value = 3017
$ws_ALC.Range("I23").Value = 3017
$ws_ALC.Range("J23").Value = 0
$ws_ALC.Range("K23").Value = 3017
$ws_ALC.Range("L23").Value = 0
$ws_ALC.Range("M23").Value = -2783

# Row 64 (ALC), Leve Item ID 5506
$ws_ALC.Range("H64").Value = 13075.154
$ws_ALC.Range("I64").Value = 7597.6
$ws_ALC.Range("J64").Value = 16498.625
$ws_ALC.Range("K64").Value = 7597.6
$ws_ALC.Range("L64").Value = 16498.625
$ws_ALC.Range("M64").Value = -7349.6
$ws_ALC.Range("N64").Value = -16994.625

# Row 67 (ALC), Leve Item ID 5506
$ws_ALC.Range("H67").Value = 13075.154
$ws_ALC.Range("I67").Value = 7597.6
$ws_ALC.Range("J67").Value = 16498.625
$ws_ALC.Range("K67").Value = 7597.6
$ws_ALC.Range("L67").Value = 16498.625
$ws_ALC.Range("M67").Value = -6739.6
$ws_ALC.Range("N67").Value = -18214.625

# Row 2 (ARM), Leve Item ID 27713
$ws_ARM.Range("H2").Value = 2288.4
$ws_ARM.Range("I2").Value = 139
$ws_ARM.Range("J2").Value = 5512.5
$ws_ARM.Range("K2").Value = 139
$ws_ARM.Range("L2").Value = 5512.5
$ws_ARM.Range("M2").Value = -26
$ws_ARM.Range("N2").Value = -5738.5

# Row 3 (ARM), Leve Item ID 2494
$ws_ARM.Range("H3").Value = 1075
$ws_ARM.Range("I3").Value = 593.75
$ws_ARM.Range("J3").Value = 3000
$ws_ARM.Range("K3").Value = 593.75
$ws_ARM.Range("L3").Value = 3000
$ws_ARM.Range("M3").Value = -478.75
$ws_ARM.Range("N3").Value = -3230

# Row 4 (ARM), Leve Item ID 5071
$ws_ARM.Range("H4").Value = 374
$ws_ARM.Range("I4").Value = 19.5
$ws_ARM.Range("J4").Value = 728.5
$ws_ARM.Range("K4").Value = 19.5
$ws_ARM.Range("L4").Value = 728.5
$ws_ARM.Range("M4").Value = 96.5
$ws_ARM.Range("N4").Value = -960.5

# Row 6 (ARM), Leve Item ID 2226
$ws_ARM.Range("H6").Value = 3430
$ws_ARM.Range("I6").Value = 290
$ws_ARM.Range("J6").Value = 5000
$ws_ARM.Range("K6").Value = 290
$ws_ARM.Range("L6").Value = 5000
$ws_ARM.Range("M6").Value = -117
$ws_ARM.Range("N6").Value = -5346

# Row 36 (ARM), Leve Item ID 3068
$ws_ARM.Range("H36").Value = 18997.5
$ws_ARM.Range("I36").Value = 18997.5
$ws_ARM.Range("J36").Value = 0
$ws_ARM.Range("K36").Value = 18997.5
$ws_ARM.Range("L36").Value = 0
$ws_ARM.Range("M36").Value = -18651.5

# Row 45 (ARM), Leve Item ID 27714
$ws_ARM.Range("H45").Value = 6214.2856
$ws_ARM.Range("I45").Value = 1500
$ws_ARM.Range("J45").Value = 7000
$ws_ARM.Range("K45").Value = 1500
$ws_ARM.Range("L45").Value = 7000
$ws_ARM.Range("M45").Value = -1123
$ws_ARM.Range("N45").Value = -7754

# Row 46 (ARM), Leve Item ID 3498
$ws_ARM.Range("H46").Value = 21316.666
$ws_ARM.Range("I46").Value = 0
$ws_ARM.Range("J46").Value = 21316.666
$ws_ARM.Range("K46").Value = 0
$ws_ARM.Range("L46").Value = 21316.666
$ws_ARM.Range("N46").Value = -21954.666
$ws_ARM.Range("M46").ClearContents()

# Row 101 (ARM), Leve Item ID 18518
$ws_ARM.Range("H101").Value = 0
$ws_ARM.Range("I101").Value = 0
$ws_ARM.Range("J101").Value = 0
$ws_ARM.Range("K101").Value = 0
$ws_ARM.Range("L101").Value = 0
$ws_ARM.Range("N101").ClearContents()

# Row 116 (ARM), Leve Item ID 27713
$ws_ARM.Range("H116").Value = 2288.4
$ws_ARM.Range("I116").Value = 139
$ws_ARM.Range("J116").Value = 5512.5
$ws_ARM.Range("K116").Value = 139
$ws_ARM.Range("L116").Value = 5512.5
$ws_ARM.Range("M116").Value = 2155
$ws_ARM.Range("N116").Value = -10100.5

# Row 122 (ARM), Leve Item ID 36168
$ws_ARM.Range("H122").Value = 3001
$ws_ARM.Range("I122").Value = 3001
$ws_ARM.Range("J122").Value = 0
$ws_ARM.Range("K122").Value = 9003
$ws_ARM.Range("L122").Value = 0
$ws_ARM.Range("M122").Value = -6553

# Row 132 (ARM), Leve Item ID 43997
$ws_ARM.Range("H132").Value = 2735.889
$ws_ARM.Range("I132").Value = 2944.1428
$ws_ARM.Range("J132").Value = 2007
$ws_ARM.Range("K132").Value = 8832.428400000001
$ws_ARM.Range("L132").Value = 6021
$ws_ARM.Range("M132").Value = -6302.428400000001
$ws_ARM.Range("N132").Value = -11081

# Row 3 (BSM), Leve Item ID 27713
$ws_BSM.Range("H3").Value = 2288.4
$ws_BSM.Range("I3").Value = 139
$ws_BSM.Range("J3").Value = 5512.5
$ws_BSM.Range("K3").Value = 139
$ws_BSM.Range("L3").Value = 5512.5
$ws_BSM.Range("M3").Value = -25
$ws_BSM.Range("N3").Value = -5740.5

# Row 10 (BSM), Leve Item ID 2417
$ws_BSM.Range("H10").Value = 608
$ws_BSM.Range("I10").Value = 220
$ws_BSM.Range("J10").Value = 996
$ws_BSM.Range("K10").Value = 220
$ws_BSM.Range("L10").Value = 996
$ws_BSM.Range("M10").Value = -80
$ws_BSM.Range("N10").Value = -1276

# Row 12 (BSM), Leve Item ID 2392
$ws_BSM.Range("H12").Value = 916.3333
$ws_BSM.Range("I12").Value = 474.5
$ws_BSM.Range("J12").Value = 1137.25
$ws_BSM.Range("K12").Value = 474.5
$ws_BSM.Range("L12").Value = 1137.25
$ws_BSM.Range("M12").Value = -306.5
$ws_BSM.Range("N12").Value = -1473.25

# Row 23 (BSM), Leve Item ID 1686
$ws_BSM.Range("H23").Value = 0
$ws_BSM.Range("I23").Value = 0
$ws_BSM.Range("J23").Value = 0
$ws_BSM.Range("K23").Value = 0
$ws_BSM.Range("L23").Value = 0
$ws_BSM.Range("N23").ClearContents()

# Row 26 (BSM), Leve Item ID 19535
$ws_BSM.Range("H26").Value = 30000
$ws_BSM.Range("I26").Value = 30000
$ws_BSM.Range("J26").Value = 0
$ws_BSM.Range("K26").Value = 30000
$ws_BSM.Range("L26").Value = 0
$ws_BSM.Range("M26").Value = -29708

# Row 96 (BSM), Leve Item ID 19525
$ws_BSM.Range("H96").Value = 10499.667
$ws_BSM.Range("I96").Value = 10499.667
$ws_BSM.Range("J96").Value = 0
$ws_BSM.Range("K96").Value = 10499.667
$ws_BSM.Range("L96").Value = 0
$ws_BSM.Range("M96").Value = -7753.666999999999
$ws_BSM.Range("N96").ClearContents()

# Row 107 (BSM), Leve Item ID 27706
$ws_BSM.Range("H107").Value = 4366.3335
$ws_BSM.Range("I107").Value = 4366.3335
$ws_BSM.Range("J107").Value = 0
$ws_BSM.Range("K107").Value = 4366.3335
$ws_BSM.Range("L107").Value = 0
$ws_BSM.Range("M107").Value = -2446.3335

# Row 134 (BSM), Leve Item ID 43998
$ws_BSM.Range("H134").Value = 4697.4287
$ws_BSM.Range("I134").Value = 4577.6
$ws_BSM.Range("J134").Value = 4997
$ws_BSM.Range("K134").Value = 13732.8
$ws_BSM.Range("L134").Value = 14991
$ws_BSM.Range("M134").Value = -11197.8
$ws_BSM.Range("N134").Value = -20061

# Row 3 (CRP), Leve Item ID 3763
$ws_CRP.Range("H3").Value = 1366.5
$ws_CRP.Range("I3").Value = 799.75
$ws_CRP.Range("J3").Value = 2500
$ws_CRP.Range("K3").Value = 799.75
$ws_CRP.Range("L3").Value = 2500
$ws_CRP.Range("M3").Value = -686.75
$ws_CRP.Range("N3").Value = -2726

# Row 12 (CUL), Leve Item ID 4854
$ws_CUL.Range("H12").Value = 1082.8334
$ws_CUL.Range("I12").Value = 41
$ws_CUL.Range("J12").Value = 1177.5454
$ws_CUL.Range("K12").Value = 123
$ws_CUL.Range("L12").Value = 3532.6362
$ws_CUL.Range("M12").Value = 50
$ws_CUL.Range("N12").Value = -3878.6362

# Row 25 (CUL), Leve Item ID 4709
$ws_CUL.Range("H25").Value = 0
$ws_CUL.Range("I25").Value = 0
$ws_CUL.Range("J25").Value = 0
$ws_CUL.Range("K25").Value = 0
$ws_CUL.Range("L25").Value = 0
$ws_CUL.Range("M25").ClearContents()
$ws_CUL.Range("N25").ClearContents()

# Row 30 (CUL), Leve Item ID 4709
$ws_CUL.Range("H30").Value = 0
$ws_CUL.Range("I30").Value = 0
$ws_CUL.Range("J30").Value = 0
$ws_CUL.Range("K30").Value = 0
$ws_CUL.Range("L30").Value = 0
$ws_CUL.Range("M30").ClearContents()
$ws_CUL.Range("N30").ClearContents()

# Row 33 (CUL), Leve Item ID 4867
$ws_CUL.Range("H33").Value = 224.25
$ws_CUL.Range("I33").Value = 199
$ws_CUL.Range("J33").Value = 249.5
$ws_CUL.Range("K33").Value = 1194
$ws_CUL.Range("L33").Value = 1497
$ws_CUL.Range("M33").Value = -911
$ws_CUL.Range("N33").Value = -2063

# Row 35 (CUL), Leve Item ID 4718
$ws_CUL.Range("H35").Value = 3613.3333
$ws_CUL.Range("I35").Value = 1870
$ws_CUL.Range("J35").Value = 3831.25
$ws_CUL.Range("K35").Value = 5610
$ws_CUL.Range("L35").Value = 11493.75
$ws_CUL.Range("M35").Value = -5322
$ws_CUL.Range("N35").Value = -12069.75

# Row 38 (CUL), Leve Item ID 4860
$ws_CUL.Range("H38").Value = 328.14285
$ws_CUL.Range("I38").Value = 99
$ws_CUL.Range("J38").Value = 419.8
$ws_CUL.Range("K38").Value = 297
$ws_CUL.Range("L38").Value = 1259.4
$ws_CUL.Range("M38").Value = 50
$ws_CUL.Range("N38").Value = -1953.4

# Row 46 (CUL), Leve Item ID 4701
$ws_CUL.Range("H46").Value = 0
$ws_CUL.Range("I46").Value = 0
$ws_CUL.Range("J46").Value = 0
$ws_CUL.Range("K46").Value = 0
$ws_CUL.Range("L46").Value = 0
$ws_CUL.Range("M46").ClearContents()

# Row 55 (CUL), Leve Item ID 4733
$ws_CUL.Range("H55").Value = 1309
$ws_CUL.Range("I55").Value = 399
$ws_CUL.Range("J55").Value = 1400
$ws_CUL.Range("K55").Value = 1197
$ws_CUL.Range("L55").Value = 4200
$ws_CUL.Range("M55").Value = -1020
$ws_CUL.Range("N55").Value = -4554

# Row 60 (CUL), Leve Item ID 4750
$ws_CUL.Range("H60").Value = 150
$ws_CUL.Range("I60").Value = 150
$ws_CUL.Range("J60").Value = 0
$ws_CUL.Range("K60").Value = 450
$ws_CUL.Range("L60").Value = 0
$ws_CUL.Range("M60").Value = -199

# Row 97 (CUL), Leve Item ID 19846
$ws_CUL.Range("H97").Value = 967.6
$ws_CUL.Range("I97").Value = 981.6667
$ws_CUL.Range("J97").Value = 946.5
$ws_CUL.Range("K97").Value = 2945.0001
$ws_CUL.Range("L97").Value = 2839.5
$ws_CUL.Range("M97").Value = -2449.0001
$ws_CUL.Range("N97").Value = -3831.5

# Row 132 (CUL), Leve Item ID 43972
$ws_CUL.Range("H132").Value = 2134.2856
$ws_CUL.Range("I132").Value = 1197
$ws_CUL.Range("J132").Value = 2509.2
$ws_CUL.Range("K132").Value = 10773
$ws_CUL.Range("L132").Value = 22582.8
$ws_CUL.Range("M132").Value = -8243
$ws_CUL.Range("N132").Value = -27642.8

# Row 6 (GSM), Leve Item ID 2108
$ws_GSM.Range("H6").Value = 1503.25
$ws_GSM.Range("I6").Value = 1336
$ws_GSM.Range("J6").Value = 2005
$ws_GSM.Range("K6").Value = 1336
$ws_GSM.Range("L6").Value = 2005
$ws_GSM.Range("M6").Value = -1223
$ws_GSM.Range("N6").Value = -2231

# Row 13 (GSM), Leve Item ID 2443
$ws_GSM.Range("H13").Value = 217.2
$ws_GSM.Range("I13").Value = 120.333336
$ws_GSM.Range("J13").Value = 362.5
$ws_GSM.Range("K13").Value = 120.333336
$ws_GSM.Range("L13").Value = 362.5
$ws_GSM.Range("M13").Value = 18.666664
$ws_GSM.Range("N13").Value = -640.5

# Row 16 (GSM), Leve Item ID 2108
$ws_GSM.Range("H16").Value = 1503.25
$ws_GSM.Range("I16").Value = 1336
$ws_GSM.Range("J16").Value = 2005
$ws_GSM.Range("K16").Value = 1336
$ws_GSM.Range("L16").Value = 2005
$ws_GSM.Range("M16").Value = -1086
$ws_GSM.Range("N16").Value = -2505

# Row 102 (GSM), Leve Item ID 36169
$ws_GSM.Range("H102").Value = 0
$ws_GSM.Range("I102").Value = 0
$ws_GSM.Range("J102").Value = 0
$ws_GSM.Range("K102").Value = 0
$ws_GSM.Range("L102").Value = 0
$ws_GSM.Range("M102").ClearContents()

# Row 123 (GSM), Leve Item ID 34150
$ws_GSM.Range("H123").Value = 39999
$ws_GSM.Range("I123").Value = 0
$ws_GSM.Range("J123").Value = 39999
$ws_GSM.Range("K123").Value = 0
$ws_GSM.Range("L123").Value = 39999
$ws_GSM.Range("N123").Value = -44899

# Row 132 (GSM), Leve Item ID 44008
$ws_GSM.Range("H132").Value = 3441.3333
$ws_GSM.Range("I132").Value = 3411.4707
$ws_GSM.Range("J132").Value = 3949
$ws_GSM.Range("K132").Value = 10234.4121
$ws_GSM.Range("L132").Value = 11847
$ws_GSM.Range("M132").Value = -7704.4121
$ws_GSM.Range("N132").Value = -16907

# Row 46 (LTW), Leve Item ID 5282
$ws_LTW.Range("H46").Value = 2761.4614
$ws_LTW.Range("I46").Value = 2100
$ws_LTW.Range("J46").Value = 3055.4443
$ws_LTW.Range("K46").Value = 2100
$ws_LTW.Range("L46").Value = 3055.4443
$ws_LTW.Range("M46").Value = -1912
$ws_LTW.Range("N46").Value = -3431.4443

# Row 55 (LTW), Leve Item ID 5284
$ws_LTW.Range("H55").Value = 1122.3334
$ws_LTW.Range("I55").Value = 0
$ws_LTW.Range("J55").Value = 1122.3334
$ws_LTW.Range("K55").Value = 0
$ws_LTW.Range("L55").Value = 1122.3334
$ws_LTW.Range("N55").Value = -1468.3334

# Row 132 (LTW), Leve Item ID 44058
$ws_LTW.Range("H132").Value = 3176.75
$ws_LTW.Range("I132").Value = 3059.2144
$ws_LTW.Range("J132").Value = 3999.5
$ws_LTW.Range("K132").Value = 9177.643199999999
$ws_LTW.Range("L132").Value = 11998.5
$ws_LTW.Range("M132").Value = -6647.643199999999
$ws_LTW.Range("N132").Value = -17058.5

# Row 4 (WVR), Leve Item ID 2996
$ws_WVR.Range("H4").Value = 100
$ws_WVR.Range("I4").Value = 100
$ws_WVR.Range("J4").Value = 0
$ws_WVR.Range("K4").Value = 100
$ws_WVR.Range("L4").Value = 0
$ws_WVR.Range("M4").Value = 13

# Row 23 (WVR), Leve Item ID 3325
$ws_WVR.Range("H23").Value = 318.4
$ws_WVR.Range("I23").Value = 286.75
$ws_WVR.Range("J23").Value = 445
$ws_WVR.Range("K23").Value = 286.75
$ws_WVR.Range("L23").Value = 445
$ws_WVR.Range("M23").Value = -57.75
$ws_WVR.Range("N23").Value = -903
